# Scheduled refresh of market-price-derived figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the per-job leve-profit sheets.
# Pure data refresh: no formulas, formatting, or structure involved.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1000.4839
$ws.Range("I107").Value = 1143.8334
$ws.Range("J107").Value = 802
$ws.Range("K107").Value = 1143.8334
$ws.Range("L107").Value = 802
$ws.Range("M107").Value = 776.1666
$ws.Range("N107").Value = -4642
$ws.Range("H111").Value = 725
$ws.Range("I111").Value = 543.5
$ws.Range("J111").Value = 1088
$ws.Range("K111").Value = 1630.5
$ws.Range("L111").Value = 3264
$ws.Range("M111").Value = 1436.5
$ws.Range("N111").Value = -9398
$ws.Range("H112").Value = 52632760
$ws.Range("J112").Value = 66668000
$ws.Range("L112").Value = 200004000
$ws.Range("N112").Value = -200006216
$ws.Range("H116").Value = 2655.1724
$ws.Range("J116").Value = 2750
$ws.Range("L116").Value = 2750
$ws.Range("N116").Value = -9634
$ws.Range("H129").Value = 1536.68
$ws.Range("J129").Value = 1559.0416
$ws.Range("L129").Value = 4677.1248
$ws.Range("N129").Value = -14677.1248
$ws.Range("H137").Value = 2037.44
$ws.Range("I137").Value = 1213.3125
$ws.Range("J137").Value = 3502.5557
$ws.Range("K137").Value = 3639.9375
$ws.Range("L137").Value = 10507.6671
$ws.Range("M137").Value = -1089.9375
$ws.Range("N137").Value = -15607.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1635644.6
$ws.Range("I2").Value = 2536.375
$ws.Range("J2").Value = 2942131.2
$ws.Range("K2").Value = 2536.375
$ws.Range("L2").Value = 2942131.2
$ws.Range("M2").Value = -2423.375
$ws.Range("N2").Value = -2942357.2
$ws.Range("H32").Value = 5412.852
$ws.Range("I32").Value = 6498.054
$ws.Range("K32").Value = 6498.054
$ws.Range("M32").Value = -6211.054
$ws.Range("H45").Value = 63418.062
$ws.Range("I45").Value = 125762.75
$ws.Range("J45").Value = 1073.375
$ws.Range("K45").Value = 125762.75
$ws.Range("L45").Value = 1073.375
$ws.Range("M45").Value = -125385.75
$ws.Range("N45").Value = -1827.375
$ws.Range("H74").Value = 4042.7896
$ws.Range("I74").Value = 1304
$ws.Range("J74").Value = 4661.2256
$ws.Range("K74").Value = 1304
$ws.Range("L74").Value = 4661.2256
$ws.Range("M74").Value = -430
$ws.Range("N74").Value = -6409.2256
$ws.Range("H77").Value = 4042.7896
$ws.Range("I77").Value = 1304
$ws.Range("J77").Value = 4661.2256
$ws.Range("K77").Value = 6520
$ws.Range("L77").Value = 23306.128
$ws.Range("M77").Value = -2152
$ws.Range("N77").Value = -32042.128
$ws.Range("H110").Value = 1050.15
$ws.Range("I110").Value = 933.7646999999999
$ws.Range("J110").Value = 1709.6666
$ws.Range("K110").Value = 933.7646999999999
$ws.Range("L110").Value = 1709.6666
$ws.Range("M110").Value = 1111.2353
$ws.Range("N110").Value = -5799.6666
$ws.Range("H116").Value = 1635644.6
$ws.Range("I116").Value = 2536.375
$ws.Range("J116").Value = 2942131.2
$ws.Range("K116").Value = 2536.375
$ws.Range("L116").Value = 2942131.2
$ws.Range("M116").Value = -242.375
$ws.Range("N116").Value = -2946719.2
$ws.Range("H130").Value = 22209.334
$ws.Range("J130").Value = 22209.334
$ws.Range("L130").Value = 22209.334
$ws.Range("N130").Value = -32249.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1635644.6
$ws.Range("I3").Value = 2536.375
$ws.Range("J3").Value = 2942131.2
$ws.Range("K3").Value = 2536.375
$ws.Range("L3").Value = 2942131.2
$ws.Range("M3").Value = -2422.375
$ws.Range("N3").Value = -2942359.2
$ws.Range("H21").Value = 30000
$ws.Range("J21").Value = 30000
$ws.Range("L21").Value = 30000
$ws.Range("N21").Value = -30472
$ws.Range("H134").Value = 49542
$ws.Range("I134").Value = 53719.566
$ws.Range("K134").Value = 161158.698
$ws.Range("M134").Value = -158623.698

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1227.5555
$ws.Range("I99").Value = 1178.2858
$ws.Range("J99").Value = 1400
$ws.Range("K99").Value = 1178.2858
$ws.Range("L99").Value = 1400
$ws.Range("M99").Value = 319.7141999999999
$ws.Range("N99").Value = -4396
$ws.Range("H107").Value = 429.3913
$ws.Range("I107").Value = 407.27777
$ws.Range("K107").Value = 407.27777
$ws.Range("M107").Value = 1512.72223
$ws.Range("H126").Value = 1227.5555
$ws.Range("I126").Value = 1178.2858
$ws.Range("J126").Value = 1400
$ws.Range("K126").Value = 3534.8574
$ws.Range("L126").Value = 4200
$ws.Range("M126").Value = -1064.8574
$ws.Range("N126").Value = -9140
$ws.Range("H130").Value = 38935
$ws.Range("J130").Value = 38935
$ws.Range("L130").Value = 38935
$ws.Range("N130").Value = -48975
$ws.Range("H132").Value = 3551.5908
$ws.Range("I132").Value = 3382.8667
$ws.Range("J132").Value = 3913.1428
$ws.Range("K132").Value = 10148.6001
$ws.Range("L132").Value = 11739.4284
$ws.Range("M132").Value = -7618.6001
$ws.Range("N132").Value = -16799.4284
$ws.Range("H134").Value = 3239.5417
$ws.Range("I134").Value = 3276.0435
$ws.Range("K134").Value = 9828.130500000001
$ws.Range("M134").Value = -7293.130500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1189.375
$ws.Range("I5").Value = 1310
$ws.Range("J5").Value = 948.125
$ws.Range("K5").Value = 3930
$ws.Range("L5").Value = 2844.375
$ws.Range("M5").Value = -3818
$ws.Range("N5").Value = -3068.375
$ws.Range("H49").Value = 2950
$ws.Range("J49").Value = 2950
$ws.Range("L49").Value = 8850
$ws.Range("N49").Value = -9162
$ws.Range("H92").Value = 170.95238
$ws.Range("I92").Value = 143.5625
$ws.Range("J92").Value = 258.6
$ws.Range("K92").Value = 430.6875
$ws.Range("L92").Value = 775.8000000000001
$ws.Range("M92").Value = 817.3125
$ws.Range("N92").Value = -3271.8
$ws.Range("H105").Value = 908000000
$ws.Range("J105").Value = 908000000
$ws.Range("L105").Value = 2724000000
$ws.Range("N105").Value = -2724005242
$ws.Range("H110").Value = 2183.3333
$ws.Range("J110").Value = 3400
$ws.Range("L110").Value = 10200
$ws.Range("N110").Value = -18380
$ws.Range("H122").Value = 11905914
$ws.Range("I122").Value = 17544336
$ws.Range("J122").Value = 2579.3333
$ws.Range("K122").Value = 157899024
$ws.Range("L122").Value = 23213.9997
$ws.Range("M122").Value = -157896574
$ws.Range("N122").Value = -28113.9997
$ws.Range("H131").Value = 1854707.5
$ws.Range("J131").Value = 2085258.9
$ws.Range("L131").Value = 6255776.699999999
$ws.Range("N131").Value = -6265856.699999999
$ws.Range("H135").Value = 1189.375
$ws.Range("I135").Value = 1310
$ws.Range("J135").Value = 948.125
$ws.Range("K135").Value = 11790
$ws.Range("L135").Value = 8533.125
$ws.Range("M135").Value = -9255
$ws.Range("N135").Value = -13603.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1448.1875
$ws.Range("I97").Value = 1513.3334
$ws.Range("J97").Value = 1252.75
$ws.Range("K97").Value = 1513.3334
$ws.Range("L97").Value = 1252.75
$ws.Range("M97").Value = -1017.3334
$ws.Range("N97").Value = -2244.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H41").Value = 13666.667
$ws.Range("J41").Value = 13000
$ws.Range("L41").Value = 13000
$ws.Range("N41").Value = -13876
$ws.Range("H46").Value = 1109.091
$ws.Range("I46").Value = 1428.5714
$ws.Range("J46").Value = 550
$ws.Range("K46").Value = 1428.5714
$ws.Range("L46").Value = 550
$ws.Range("M46").Value = -1240.5714
$ws.Range("N46").Value = -926
$ws.Range("H61").Value = 714.36
$ws.Range("I61").Value = 614.1111
$ws.Range("J61").Value = 972.1429000000001
$ws.Range("K61").Value = 614.1111
$ws.Range("L61").Value = 972.1429000000001
$ws.Range("M61").Value = -412.1111
$ws.Range("N61").Value = -1376.1429
$ws.Range("H93").Value = 3026
$ws.Range("I93").Value = 2000
$ws.Range("J93").Value = 3710
$ws.Range("K93").Value = 2000
$ws.Range("L93").Value = 3710
$ws.Range("M93").Value = -752
$ws.Range("N93").Value = -6206
$ws.Range("H100").Value = 1193.2142
$ws.Range("I100").Value = 1244.5555
$ws.Range("J100").Value = 1100.8
$ws.Range("K100").Value = 1244.5555
$ws.Range("L100").Value = 1100.8
$ws.Range("M100").Value = -703.5554999999999
$ws.Range("N100").Value = -2182.8
$ws.Range("H113").Value = 714.36
$ws.Range("I113").Value = 614.1111
$ws.Range("J113").Value = 972.1429000000001
$ws.Range("K113").Value = 614.1111
$ws.Range("L113").Value = 972.1429000000001
$ws.Range("M113").Value = 1555.8889
$ws.Range("N113").Value = -5312.1429
$ws.Range("H128").Value = 44359.8
$ws.Range("J128").Value = 44359.8
$ws.Range("L128").Value = 44359.8
$ws.Range("N128").Value = -54319.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 442.23077
$ws.Range("I100").Value = 374.83334
$ws.Range("K100").Value = 749.66668
$ws.Range("M100").Value = -208.66668
$ws.Range("H104").Value = 24533.5
$ws.Range("J104").Value = 24533.5
$ws.Range("L104").Value = 24533.5
$ws.Range("N104").Value = -31521.5
